$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4161528050899506
$ws.Range("B1").Value = 0.7137250304222107
$ws.Range("C1").Value = 0.9672904014587402
$ws.Range("D1").Value = 4.545559883117676
$ws.Range("E1").Value = 1.288611054420471
